# Update column F (dSF) values for a set of rows, per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    3  = 0
    15 = 2
    17 = 1
    19 = 2
    21 = 1
    30 = 0
    33 = 2
    41 = 1
    43 = 5
    44 = 0
    47 = 0
    53 = 0
    59 = 2
    64 = 2
    65 = 1
    70 = 3
    74 = 9
    78 = -2
    79 = 11
    80 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
